# Apply cryptos list update (cell value changes) via Excel COM interop
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.597.72"
$ws.Range("E2").Value = "  +2.19%  "
$ws.Range("D3").Value = "3.103.51"
$ws.Range("E3").Value = "  +1.03%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'528.25"
$ws.Range("E5").Value = "  +2.78%  "
$ws.Range("D6").Value = "'143.41"
$ws.Range("E6").Value = "  +2.07%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +1.48%  "
$ws.Range("E10").Value = "  +1.18%  "
$ws.Range("D11").Value = "'0.384"
$ws.Range("E11").Value = "  +3.24%  "
$ws.Range("D12").Value = "3.638.37"
$ws.Range("E12").Value = "  +0.93%  "
$ws.Range("E13").Value = "  +1.01%  "
$ws.Range("D14").Value = "'26.84"
$ws.Range("E14").Value = "  +5.42%  "
$ws.Range("D15").Value = "'0.0000166"
$ws.Range("E15").Value = "  +2.58%  "
$ws.Range("D16").Value = "58.645.50"
$ws.Range("E16").Value = "  +2.12%  "
$ws.Range("D17").Value = "3.091.11"
$ws.Range("E17").Value = "  +0.78%  "
$ws.Range("E18").Value = "  +0.81%  "
$ws.Range("D19").Value = "'12.95"
$ws.Range("E19").Value = "  -0.86%  "
$ws.Range("D20").Value = "'8.09"
$ws.Range("E20").Value = "  -0.05%  "
$ws.Range("D21").Value = "'342.34"
$ws.Range("E21").Value = "  +2.92%  "
$ws.Range("D22").Value = "'1.00"
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("E23").Value = "  +1.38%  "
$ws.Range("D24").Value = "'66.05"
$ws.Range("E24").Value = "  +0.34%  "
$ws.Range("E25").Value = "  +0.81%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("D27").Value = "0.0₃0918"
$ws.Range("E27").Value = "  +1.87%  "
$ws.Range("D28").Value = "'6.73"
$ws.Range("E28").Value = "  +5.56%  "
$ws.Range("D29").Value = "'7.24"
$ws.Range("E29").Value = "  +3.01%  "
$ws.Range("E30").Value = "  +3.67%  "
$ws.Range("D31").Value = "'1.22"
$ws.Range("E31").Value = "  +4.92%  "
$ws.Range("D32").Value = "'20.98"
$ws.Range("E32").Value = "  +1.23%  "
$ws.Range("D33").Value = "'154.33"
$ws.Range("E33").Value = "  -0.36%  "
$ws.Range("E34").Value = "  +3.51%  "
$ws.Range("D35").Value = "'6.07"
$ws.Range("E35").Value = "  +3.19%  "
$ws.Range("D36").Value = "'26.98"
$ws.Range("E36").Value = "  -1.67%  "
$ws.Range("E37").Value = "  +4.77%  "
$ws.Range("D38").Value = "'0.0680"
$ws.Range("E38").Value = "  +0.99%  "
$ws.Range("D39").Value = "3.146.97"
$ws.Range("E39").Value = "  +1.16%  "
$ws.Range("B40").Value = "Mantle"
$ws.Range("C40").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D40").Value = "'0.679"
$ws.Range("E40").Value = "  +1.56%  "
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").Value = "'3.89"
$ws.Range("E41").Value = "  +1.62%  "
$ws.Range("D42").Value = "'36.90"
$ws.Range("E42").Value = "  +0.47%  "
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").Value = "'1.50"
$ws.Range("E43").Value = "  +9.03%  "
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("D45").Value = "2.301.43"
$ws.Range("E45").Value = "  +0.74%  "
$ws.Range("E46").Value = "  +1.72%  "
$ws.Range("D47").Value = "'20.98"
$ws.Range("E47").Value = "  +5.72%  "
$ws.Range("D48").Value = "'0.971"
$ws.Range("E48").Value = "  +3.68%  "
$ws.Range("D49").Value = "'6.00"
$ws.Range("E49").Value = "  +2.34%  "
$ws.Range("D50").Value = "'271.34"
$ws.Range("E50").Value = "  +9.01%  "
$ws.Range("D51").Value = "'0.752"
$ws.Range("E51").Value = "  +10.05%  "
